# Fruta / hortaliza, semanal
# Weekly update: a new week's observation is inserted as the new row 5
# (pushing the previous rows 5-13 down to 6-14), and one more historical
# observation is appended as the new last row (15).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 5 - shifts existing rows 5-13 down to 6-14.
$ws.Rows.Item(5).Insert()

# Fill in the newly inserted row 5 with this week's data.
$ws.Range("A5").Value = 10
$ws.Range("B5").Value = "Vega Modelo de Temuco"
$ws.Range("C5").Value = "La Araucanía"
$ws.Range("D5").Value = 44749
$ws.Range("D5").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E5").Value = 9
$ws.Range("F5").Value = "Fruta"
$ws.Range("G5").Value = 100108
$ws.Range("H5").Value = "Tropicales y subtropicales"
$ws.Range("I5").Value = 100108001
$ws.Range("J5").Value = "Guayaba"
$ws.Range("K5").Value = "Sin especificar"
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 120
$ws.Range("N5").Value = 2300
$ws.Range("O5").Value = 2300
$ws.Range("P5").Value = 2300
$ws.Range("Q5").Value = "$/kilo"
$ws.Range("R5").Value = "Región de Arica y Parinacota"
$ws.Range("S5").Value = 2300
$ws.Range("T5").Value = 1

# Append a new row (15) at the end with another observation.
$ws.Range("A15").Value = 10
$ws.Range("B15").Value = "Vega Modelo de Temuco"
$ws.Range("C15").Value = "La Araucanía"
$ws.Range("D15").Value = 44748
$ws.Range("D15").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E15").Value = 9
$ws.Range("F15").Value = "Fruta"
$ws.Range("G15").Value = 100108
$ws.Range("H15").Value = "Tropicales y subtropicales"
$ws.Range("I15").Value = 100108001
$ws.Range("J15").Value = "Guayaba"
$ws.Range("K15").Value = "Sin especificar"
$ws.Range("L15").Value = "Primera"
$ws.Range("M15").Value = 300
$ws.Range("N15").Value = 2300
$ws.Range("O15").Value = 2300
$ws.Range("P15").Value = 2300
$ws.Range("Q15").Value = "$/kilo"
$ws.Range("R15").Value = "Región de Arica y Parinacota"
$ws.Range("S15").Value = 2300
$ws.Range("T15").Value = 1
